$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two header cells ("Coln names")
$ws.Range("D1").Value = "Full Team Name Digital"
$ws.Range("H1").Value = "Full Team Name EMU"

# Bold the header row and give it a yellow fill.
# Build the combined style on A1 first, then copy/paste the formatting
# across the rest of the header row so only a single new cell style is
# materialised (instead of one per incremental property write).
$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.Interior.Color = 65535
$a1.Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

# Move the active selection to H1
$ws.Range("H1").Select()
